$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A4").Value = "Tlapa"
$ws.Range("B4").Value = 21
[void]$ws.Range("B5").Select()
